$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure all touched cells are formatted as Text before assigning,
# so Excel does not reinterpret numeric-looking strings (e.g. "18.53", "0.0625") as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.002.05'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.638.10'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.37%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.34%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.52'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.95%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.54%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.59%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0625'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.53'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -5.71%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.59%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.864.40'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.41%  '
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.21'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.86%  '
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.632.49'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.531'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.70%  '
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₃0746'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.53%  '
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.001.23'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.90'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.13%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '191.83'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.58%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.23%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.10%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.10%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.12%  '
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.88'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.34%  '
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.79'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.97%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.26%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.86'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.62%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.27'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.65%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.50%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.99%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.16'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -4.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.50'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.41'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.55%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.137.19'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.35%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.869'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -4.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.521'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.52%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '98.50'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.779'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.774.07'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.40%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -4.62%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '55.25'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.30%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0528'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.48%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.50'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.13%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.71%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.03%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.04%  '
